$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 430; this shifts the existing rows 430-441
# down to become rows 434-445 (values/formatting move with them).
$ws.Range("A430:T433").EntireRow.Insert()

# --- New row 430 ---
$ws.Range("A430").Value = 9
$ws.Range("B430").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C430").Value = "Metropolitana"
$ws.Range("D430").Value = 44509
$ws.Range("E430").Value = 13
$ws.Range("F430").Value = "Fruta"
$ws.Range("G430").Value = 100101
$ws.Range("H430").Value = "Berries"
$ws.Range("I430").Value = 100101007
$ws.Range("J430").Value = "Kiwi"
$ws.Range("K430").Value = "Hayward"
$ws.Range("L430").Value = "Especial"
$ws.Range("M430").Value = 410
$ws.Range("N430").Value = 11000
$ws.Range("O430").Value = 11000
$ws.Range("P430").Value = 11000
$ws.Range("Q430").Value = "$/bandeja 10 kilos"
$ws.Range("R430").Value = "Provincia de Curicó"
$ws.Range("S430").Value = 1100
$ws.Range("T430").Value = 10

# --- New row 431 ---
$ws.Range("A431").Value = 9
$ws.Range("B431").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C431").Value = "Metropolitana"
$ws.Range("D431").Value = 44509
$ws.Range("E431").Value = 13
$ws.Range("F431").Value = "Fruta"
$ws.Range("G431").Value = 100101
$ws.Range("H431").Value = "Berries"
$ws.Range("I431").Value = 100101007
$ws.Range("J431").Value = "Kiwi"
$ws.Range("K431").Value = "Hayward"
$ws.Range("L431").Value = "Extra (doble especial)"
$ws.Range("M431").Value = 380
$ws.Range("N431").Value = 12000
$ws.Range("O431").Value = 12000
$ws.Range("P431").Value = 12000
$ws.Range("Q431").Value = "$/bandeja 10 kilos"
$ws.Range("R431").Value = "Provincia de Curicó"
$ws.Range("S431").Value = 1200
$ws.Range("T431").Value = 10

# --- New row 432 ---
$ws.Range("A432").Value = 9
$ws.Range("B432").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C432").Value = "Metropolitana"
$ws.Range("D432").Value = 44509
$ws.Range("E432").Value = 13
$ws.Range("F432").Value = "Fruta"
$ws.Range("G432").Value = 100101
$ws.Range("H432").Value = "Berries"
$ws.Range("I432").Value = 100101007
$ws.Range("J432").Value = "Kiwi"
$ws.Range("K432").Value = "Hayward"
$ws.Range("L432").Value = "Primera"
$ws.Range("M432").Value = 440
$ws.Range("N432").Value = 10000
$ws.Range("O432").Value = 10000
$ws.Range("P432").Value = 10000
$ws.Range("Q432").Value = "$/bandeja 10 kilos"
$ws.Range("R432").Value = "Provincia de Curicó"
$ws.Range("S432").Value = 1000
$ws.Range("T432").Value = 10

# --- New row 433 ---
$ws.Range("A433").Value = 9
$ws.Range("B433").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C433").Value = "Metropolitana"
$ws.Range("D433").Value = 44509
$ws.Range("E433").Value = 13
$ws.Range("F433").Value = "Fruta"
$ws.Range("G433").Value = 100101
$ws.Range("H433").Value = "Berries"
$ws.Range("I433").Value = 100101007
$ws.Range("J433").Value = "Kiwi"
$ws.Range("K433").Value = "Hayward"
$ws.Range("L433").Value = "Segunda"
$ws.Range("M433").Value = 350
$ws.Range("N433").Value = 7000
$ws.Range("O433").Value = 7000
$ws.Range("P433").Value = 7000
$ws.Range("Q433").Value = "$/bandeja 10 kilos"
$ws.Range("R433").Value = "Provincia de Curicó"
$ws.Range("S433").Value = 700
$ws.Range("T433").Value = 10
